$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the F:V contents of rows 15 and 16 (home/away match data was mismatched) ---
$ws.Cells.Item(15, 6).Value = "Buriram"
$ws.Cells.Item(15, 7).Value = 3
$ws.Cells.Item(15, 8).Value = "Lamphun Warrior"
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 1.25
$ws.Cells.Item(15, 11).Value = "15/08/2023 16:42"
$ws.Cells.Item(15, 12).Value = 1.21
$ws.Cells.Item(15, 13).Value = "20/08/2023 12:50"
$ws.Cells.Item(15, 14).Value = 6.02
$ws.Cells.Item(15, 15).Value = "15/08/2023 16:42"
$ws.Cells.Item(15, 16).Value = 6.56
$ws.Cells.Item(15, 17).Value = "20/08/2023 12:50"
$ws.Cells.Item(15, 18).Value = 10.71
$ws.Cells.Item(15, 19).Value = "15/08/2023 16:42"
$ws.Cells.Item(15, 20).Value = 12.97
$ws.Cells.Item(15, 21).Value = "20/08/2023 12:50"
$ws.Cells.Item(15, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/buriram-united-f-c-lamphun-warrior/nLdbew0s/"

$ws.Cells.Item(16, 6).Value = "Sukhothai"
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = "Trat FC"
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 1.78
$ws.Cells.Item(16, 11).Value = "18/08/2023 23:26"
$ws.Cells.Item(16, 12).Value = 1.83
$ws.Cells.Item(16, 13).Value = "20/08/2023 12:51"
$ws.Cells.Item(16, 14).Value = 3.66
$ws.Cells.Item(16, 15).Value = "18/08/2023 23:26"
$ws.Cells.Item(16, 16).Value = 3.87
$ws.Cells.Item(16, 17).Value = "20/08/2023 12:51"
$ws.Cells.Item(16, 18).Value = 4.14
$ws.Cells.Item(16, 19).Value = "18/08/2023 23:26"
$ws.Cells.Item(16, 20).Value = 4.1
$ws.Cells.Item(16, 21).Value = "20/08/2023 12:51"
$ws.Cells.Item(16, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/sukhothai-trat-fc/xGYJjFEC/"

# --- Append new rows 91:109 (new match results/odds data) ---
# Copy formatting (styles) from the last existing row (90) down through row 109
$ws.Range("A90:V90").Copy() | Out-Null
$ws.Range("A91:V109").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 91
$ws.Cells.Item(91, 1).Value = 90
$ws.Cells.Item(91, 2).Value = "thailand"
$ws.Cells.Item(91, 3).Value = "thai-league-1"
$ws.Cells.Item(91, 4).Value = "2023-2024"
$ws.Cells.Item(91, 5).Value = 45263.47916666666
$ws.Cells.Item(91, 6).Value = "Khonkaen Utd."
$ws.Cells.Item(91, 7).Value = 2
$ws.Cells.Item(91, 8).Value = "Lamphun Warrior"
$ws.Cells.Item(91, 9).Value = 2
$ws.Cells.Item(91, 10).Value = 3.16
$ws.Cells.Item(91, 11).Value = "27/11/2023 12:42"
$ws.Cells.Item(91, 12).Value = 3.11
$ws.Cells.Item(91, 13).Value = "03/12/2023 11:22"
$ws.Cells.Item(91, 14).Value = 3.56
$ws.Cells.Item(91, 15).Value = "27/11/2023 12:42"
$ws.Cells.Item(91, 16).Value = 3.59
$ws.Cells.Item(91, 17).Value = "03/12/2023 11:22"
$ws.Cells.Item(91, 18).Value = 2.1
$ws.Cells.Item(91, 19).Value = "27/11/2023 12:42"
$ws.Cells.Item(91, 20).Value = 2.25
$ws.Cells.Item(91, 21).Value = "03/12/2023 11:22"
$ws.Cells.Item(91, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/khonkaen-united-lamphun-warrior/Wphv7Uan/"

# Row 92
$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = "thailand"
$ws.Cells.Item(92, 3).Value = "thai-league-1"
$ws.Cells.Item(92, 4).Value = "2023-2024"
$ws.Cells.Item(92, 5).Value = 45263.5
$ws.Cells.Item(92, 6).Value = "Muang Thong Utd"
$ws.Cells.Item(92, 7).Value = 2
$ws.Cells.Item(92, 8).Value = "Buriram"
$ws.Cells.Item(92, 9).Value = 2
$ws.Cells.Item(92, 10).Value = 3.36
$ws.Cells.Item(92, 11).Value = "26/11/2023 12:13"
$ws.Cells.Item(92, 12).Value = 4.48
$ws.Cells.Item(92, 13).Value = "03/12/2023 11:59"
$ws.Cells.Item(92, 14).Value = 3.62
$ws.Cells.Item(92, 15).Value = "26/11/2023 12:13"
$ws.Cells.Item(92, 16).Value = 4.1
$ws.Cells.Item(92, 17).Value = "03/12/2023 11:59"
$ws.Cells.Item(92, 18).Value = 2.08
$ws.Cells.Item(92, 19).Value = "26/11/2023 12:13"
$ws.Cells.Item(92, 20).Value = 1.72
$ws.Cells.Item(92, 21).Value = "03/12/2023 11:59"
$ws.Cells.Item(92, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/muang-thong-utd-buriram-united-f-c/MZpi4Sq5/"

# Row 93
$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = "thailand"
$ws.Cells.Item(93, 3).Value = "thai-league-1"
$ws.Cells.Item(93, 4).Value = "2023-2024"
$ws.Cells.Item(93, 5).Value = 45263.54166666666
$ws.Cells.Item(93, 6).Value = "Prachuap"
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = "Bangkok Utd"
$ws.Cells.Item(93, 9).Value = 1
$ws.Cells.Item(93, 10).Value = 3.82
$ws.Cells.Item(93, 11).Value = "26/11/2023 13:12"
$ws.Cells.Item(93, 12).Value = 4.83
$ws.Cells.Item(93, 13).Value = "03/12/2023 12:52"
$ws.Cells.Item(93, 14).Value = 3.76
$ws.Cells.Item(93, 15).Value = "26/11/2023 13:12"
$ws.Cells.Item(93, 16).Value = 4.05
$ws.Cells.Item(93, 17).Value = "03/12/2023 12:52"
$ws.Cells.Item(93, 18).Value = 1.89
$ws.Cells.Item(93, 19).Value = "26/11/2023 13:12"
$ws.Cells.Item(93, 20).Value = 1.68
$ws.Cells.Item(93, 21).Value = "03/12/2023 12:52"
$ws.Cells.Item(93, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/prachuap-bangkok-utd/Q7ir6lFh/"

# Row 94
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = "thailand"
$ws.Cells.Item(94, 3).Value = "thai-league-1"
$ws.Cells.Item(94, 4).Value = "2023-2024"
$ws.Cells.Item(94, 5).Value = 45269.5
$ws.Cells.Item(94, 6).Value = "Bangkok Utd"
$ws.Cells.Item(94, 7).Value = 4
$ws.Cells.Item(94, 8).Value = "Khonkaen Utd."
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 1.19
$ws.Cells.Item(94, 11).Value = "03/12/2023 13:12"
$ws.Cells.Item(94, 12).Value = 1.21
$ws.Cells.Item(94, 13).Value = "09/12/2023 11:55"
$ws.Cells.Item(94, 14).Value = 6.77
$ws.Cells.Item(94, 15).Value = "03/12/2023 13:12"
$ws.Cells.Item(94, 16).Value = 7.01
$ws.Cells.Item(94, 17).Value = "09/12/2023 11:55"
$ws.Cells.Item(94, 18).Value = 9.85
$ws.Cells.Item(94, 19).Value = "03/12/2023 13:12"
$ws.Cells.Item(94, 20).Value = 11.95
$ws.Cells.Item(94, 21).Value = "09/12/2023 11:55"
$ws.Cells.Item(94, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/bangkok-utd-khonkaen-united/jeOZe4io/"

# Row 95
$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = "thailand"
$ws.Cells.Item(95, 3).Value = "thai-league-1"
$ws.Cells.Item(95, 4).Value = "2023-2024"
$ws.Cells.Item(95, 5).Value = 45269.54166666666
$ws.Cells.Item(95, 6).Value = "Buriram"
$ws.Cells.Item(95, 7).Value = 4
$ws.Cells.Item(95, 8).Value = "Nakhon Pathom"
$ws.Cells.Item(95, 9).Value = 1
$ws.Cells.Item(95, 10).Value = 1.15
$ws.Cells.Item(95, 11).Value = "03/12/2023 12:13"
$ws.Cells.Item(95, 12).Value = 1.21
$ws.Cells.Item(95, 13).Value = "09/12/2023 12:51"
$ws.Cells.Item(95, 14).Value = 7.34
$ws.Cells.Item(95, 15).Value = "03/12/2023 12:13"
$ws.Cells.Item(95, 16).Value = 6.9
$ws.Cells.Item(95, 17).Value = "09/12/2023 12:51"
$ws.Cells.Item(95, 18).Value = 11.1
$ws.Cells.Item(95, 19).Value = "03/12/2023 12:13"
$ws.Cells.Item(95, 20).Value = 11.89
$ws.Cells.Item(95, 21).Value = "09/12/2023 12:51"
$ws.Cells.Item(95, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/buriram-united-f-c-nakhon-pathom/ns8xGPrH/"

# Row 96
$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = "thailand"
$ws.Cells.Item(96, 3).Value = "thai-league-1"
$ws.Cells.Item(96, 4).Value = "2023-2024"
$ws.Cells.Item(96, 5).Value = 45269.58333333334
$ws.Cells.Item(96, 6).Value = "Pathum United"
$ws.Cells.Item(96, 7).Value = 2
$ws.Cells.Item(96, 8).Value = "Ratchaburi"
$ws.Cells.Item(96, 9).Value = 1
$ws.Cells.Item(96, 10).Value = 1.8
$ws.Cells.Item(96, 11).Value = "02/12/2023 14:12"
$ws.Cells.Item(96, 12).Value = 1.67
$ws.Cells.Item(96, 13).Value = "09/12/2023 13:52"
$ws.Cells.Item(96, 14).Value = 3.66
$ws.Cells.Item(96, 15).Value = "02/12/2023 14:12"
$ws.Cells.Item(96, 16).Value = 4.1
$ws.Cells.Item(96, 17).Value = "09/12/2023 13:52"
$ws.Cells.Item(96, 18).Value = 3.98
$ws.Cells.Item(96, 19).Value = "02/12/2023 14:12"
$ws.Cells.Item(96, 20).Value = 4.85
$ws.Cells.Item(96, 21).Value = "09/12/2023 13:52"
$ws.Cells.Item(96, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/pathum-united-ratchaburi/ET9UHoE4/"

# Row 97
$ws.Cells.Item(97, 1).Value = 96
$ws.Cells.Item(97, 2).Value = "thailand"
$ws.Cells.Item(97, 3).Value = "thai-league-1"
$ws.Cells.Item(97, 4).Value = "2023-2024"
$ws.Cells.Item(97, 5).Value = 45270.47916666666
$ws.Cells.Item(97, 6).Value = "Trat FC"
$ws.Cells.Item(97, 7).Value = 2
$ws.Cells.Item(97, 8).Value = "Uthai Thani"
$ws.Cells.Item(97, 9).Value = 3
$ws.Cells.Item(97, 10).Value = 2.22
$ws.Cells.Item(97, 11).Value = "03/12/2023 11:42"
$ws.Cells.Item(97, 12).Value = 2.57
$ws.Cells.Item(97, 13).Value = "10/12/2023 11:27"
$ws.Cells.Item(97, 14).Value = 3.52
$ws.Cells.Item(97, 15).Value = "03/12/2023 11:42"
$ws.Cells.Item(97, 16).Value = 3.55
$ws.Cells.Item(97, 17).Value = "10/12/2023 11:27"
$ws.Cells.Item(97, 18).Value = 2.94
$ws.Cells.Item(97, 19).Value = "03/12/2023 11:42"
$ws.Cells.Item(97, 20).Value = 2.68
$ws.Cells.Item(97, 21).Value = "10/12/2023 11:27"
$ws.Cells.Item(97, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/trat-fc-uthai-thani/0Q5QIRbb/"

# Row 98
$ws.Cells.Item(98, 1).Value = 97
$ws.Cells.Item(98, 2).Value = "thailand"
$ws.Cells.Item(98, 3).Value = "thai-league-1"
$ws.Cells.Item(98, 4).Value = "2023-2024"
$ws.Cells.Item(98, 5).Value = 45270.5
$ws.Cells.Item(98, 6).Value = "Chonburi"
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = "Port MTI FC"
$ws.Cells.Item(98, 9).Value = 2
$ws.Cells.Item(98, 10).Value = 3.11
$ws.Cells.Item(98, 11).Value = "03/12/2023 12:13"
$ws.Cells.Item(98, 12).Value = 3.26
$ws.Cells.Item(98, 13).Value = "10/12/2023 11:59"
$ws.Cells.Item(98, 14).Value = 3.62
$ws.Cells.Item(98, 15).Value = "03/12/2023 12:13"
$ws.Cells.Item(98, 16).Value = 3.79
$ws.Cells.Item(98, 17).Value = "10/12/2023 11:51"
$ws.Cells.Item(98, 18).Value = 2.09
$ws.Cells.Item(98, 19).Value = "03/12/2023 12:13"
$ws.Cells.Item(98, 20).Value = 2.11
$ws.Cells.Item(98, 21).Value = "10/12/2023 11:37"
$ws.Cells.Item(98, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/chonburi-port-mti-fc/8j9YG5TA/"

# Row 99
$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = "thailand"
$ws.Cells.Item(99, 3).Value = "thai-league-1"
$ws.Cells.Item(99, 4).Value = "2023-2024"
$ws.Cells.Item(99, 5).Value = 45270.54166666666
$ws.Cells.Item(99, 6).Value = "Lamphun Warrior"
$ws.Cells.Item(99, 7).Value = 2
$ws.Cells.Item(99, 8).Value = "Prachuap"
$ws.Cells.Item(99, 9).Value = 1
$ws.Cells.Item(99, 10).Value = 2.4
$ws.Cells.Item(99, 11).Value = "03/12/2023 13:12"
$ws.Cells.Item(99, 12).Value = 2.35
$ws.Cells.Item(99, 13).Value = "10/12/2023 12:57"
$ws.Cells.Item(99, 14).Value = 3.41
$ws.Cells.Item(99, 15).Value = "03/12/2023 13:12"
$ws.Cells.Item(99, 16).Value = 3.4
$ws.Cells.Item(99, 17).Value = "10/12/2023 12:57"
$ws.Cells.Item(99, 18).Value = 2.77
$ws.Cells.Item(99, 19).Value = "03/12/2023 13:12"
$ws.Cells.Item(99, 20).Value = 3.08
$ws.Cells.Item(99, 21).Value = "10/12/2023 12:57"
$ws.Cells.Item(99, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/lamphun-warrior-prachuap/ShPZEqbN/"

# Row 100
$ws.Cells.Item(100, 1).Value = 99
$ws.Cells.Item(100, 2).Value = "thailand"
$ws.Cells.Item(100, 3).Value = "thai-league-1"
$ws.Cells.Item(100, 4).Value = "2023-2024"
$ws.Cells.Item(100, 5).Value = 45271.58333333334
$ws.Cells.Item(100, 6).Value = "Chiangrai Utd"
$ws.Cells.Item(100, 7).Value = 1
$ws.Cells.Item(100, 8).Value = "Police Tero"
$ws.Cells.Item(100, 9).Value = 2
$ws.Cells.Item(100, 10).Value = 2.27
$ws.Cells.Item(100, 11).Value = "04/12/2023 15:42"
$ws.Cells.Item(100, 12).Value = 1.9
$ws.Cells.Item(100, 13).Value = "11/12/2023 13:55"
$ws.Cells.Item(100, 14).Value = 3.62
$ws.Cells.Item(100, 15).Value = "04/12/2023 15:42"
$ws.Cells.Item(100, 16).Value = 3.73
$ws.Cells.Item(100, 17).Value = "11/12/2023 13:55"
$ws.Cells.Item(100, 18).Value = 2.82
$ws.Cells.Item(100, 19).Value = "04/12/2023 15:42"
$ws.Cells.Item(100, 20).Value = 3.96
$ws.Cells.Item(100, 21).Value = "11/12/2023 13:55"
$ws.Cells.Item(100, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/chiangrai-utd-police-tero/I1GLJ7qh/"

# Row 101
$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = "thailand"
$ws.Cells.Item(101, 3).Value = "thai-league-1"
$ws.Cells.Item(101, 4).Value = "2023-2024"
$ws.Cells.Item(101, 5).Value = 45272.52083333334
$ws.Cells.Item(101, 6).Value = "Sukhothai"
$ws.Cells.Item(101, 7).Value = 1
$ws.Cells.Item(101, 8).Value = "Muang Thong Utd"
$ws.Cells.Item(101, 9).Value = 2
$ws.Cells.Item(101, 10).Value = 2.66
$ws.Cells.Item(101, 11).Value = "05/12/2023 19:42"
$ws.Cells.Item(101, 12).Value = 2.46
$ws.Cells.Item(101, 13).Value = "12/12/2023 12:28"
$ws.Cells.Item(101, 14).Value = 3.48
$ws.Cells.Item(101, 15).Value = "05/12/2023 19:42"
$ws.Cells.Item(101, 16).Value = 3.84
$ws.Cells.Item(101, 17).Value = "12/12/2023 12:29"
$ws.Cells.Item(101, 18).Value = 2.44
$ws.Cells.Item(101, 19).Value = "05/12/2023 19:42"
$ws.Cells.Item(101, 20).Value = 2.64
$ws.Cells.Item(101, 21).Value = "12/12/2023 12:29"
$ws.Cells.Item(101, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/sukhothai-muang-thong-utd/vBHHKmUo/"

# Row 102
$ws.Cells.Item(102, 1).Value = 101
$ws.Cells.Item(102, 2).Value = "thailand"
$ws.Cells.Item(102, 3).Value = "thai-league-1"
$ws.Cells.Item(102, 4).Value = "2023-2024"
$ws.Cells.Item(102, 5).Value = 45275.54166666666
$ws.Cells.Item(102, 6).Value = "Prachuap"
$ws.Cells.Item(102, 7).Value = 3
$ws.Cells.Item(102, 8).Value = "Chonburi"
$ws.Cells.Item(102, 9).Value = 1
$ws.Cells.Item(102, 10).Value = 2.42
$ws.Cells.Item(102, 11).Value = "11/12/2023 09:12"
$ws.Cells.Item(102, 12).Value = 2.44
$ws.Cells.Item(102, 13).Value = "15/12/2023 12:55"
$ws.Cells.Item(102, 14).Value = 3.43
$ws.Cells.Item(102, 15).Value = "11/12/2023 09:12"
$ws.Cells.Item(102, 16).Value = 3.6
$ws.Cells.Item(102, 17).Value = "15/12/2023 12:51"
$ws.Cells.Item(102, 18).Value = 2.71
$ws.Cells.Item(102, 19).Value = "11/12/2023 09:12"
$ws.Cells.Item(102, 20).Value = 2.8
$ws.Cells.Item(102, 21).Value = "15/12/2023 12:55"
$ws.Cells.Item(102, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/prachuap-chonburi/lQrbFD93/"

# Row 103
$ws.Cells.Item(103, 1).Value = 102
$ws.Cells.Item(103, 2).Value = "thailand"
$ws.Cells.Item(103, 3).Value = "thai-league-1"
$ws.Cells.Item(103, 4).Value = "2023-2024"
$ws.Cells.Item(103, 5).Value = 45276.5
$ws.Cells.Item(103, 6).Value = "Khonkaen Utd."
$ws.Cells.Item(103, 7).Value = 1
$ws.Cells.Item(103, 8).Value = "Pathum United"
$ws.Cells.Item(103, 9).Value = 1
$ws.Cells.Item(103, 10).Value = 4.97
$ws.Cells.Item(103, 11).Value = "11/12/2023 11:12"
$ws.Cells.Item(103, 12).Value = 4.72
$ws.Cells.Item(103, 13).Value = "16/12/2023 11:54"
$ws.Cells.Item(103, 14).Value = 4.26
$ws.Cells.Item(103, 15).Value = "11/12/2023 11:12"
$ws.Cells.Item(103, 16).Value = 4.36
$ws.Cells.Item(103, 17).Value = "16/12/2023 11:56"
$ws.Cells.Item(103, 18).Value = 1.55
$ws.Cells.Item(103, 19).Value = "11/12/2023 11:12"
$ws.Cells.Item(103, 20).Value = 1.65
$ws.Cells.Item(103, 21).Value = "16/12/2023 11:56"
$ws.Cells.Item(103, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/khonkaen-united-pathum-united/OOA59NLp/"

# Row 104
$ws.Cells.Item(104, 1).Value = 103
$ws.Cells.Item(104, 2).Value = "thailand"
$ws.Cells.Item(104, 3).Value = "thai-league-1"
$ws.Cells.Item(104, 4).Value = "2023-2024"
$ws.Cells.Item(104, 5).Value = 45276.52083333334
$ws.Cells.Item(104, 6).Value = "Muang Thong Utd"
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = "Ratchaburi"
$ws.Cells.Item(104, 9).Value = 1
$ws.Cells.Item(104, 10).Value = 2.1
$ws.Cells.Item(104, 11).Value = "12/12/2023 12:42"
$ws.Cells.Item(104, 12).Value = 2.14
$ws.Cells.Item(104, 13).Value = "16/12/2023 12:27"
$ws.Cells.Item(104, 14).Value = 3.57
$ws.Cells.Item(104, 15).Value = "12/12/2023 12:42"
$ws.Cells.Item(104, 16).Value = 3.95
$ws.Cells.Item(104, 17).Value = "16/12/2023 12:27"
$ws.Cells.Item(104, 18).Value = 3.16
$ws.Cells.Item(104, 19).Value = "12/12/2023 12:42"
$ws.Cells.Item(104, 20).Value = 3.08
$ws.Cells.Item(104, 21).Value = "16/12/2023 12:27"
$ws.Cells.Item(104, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/muang-thong-utd-ratchaburi/xSE98syi/"

# Row 105
$ws.Cells.Item(105, 1).Value = 104
$ws.Cells.Item(105, 2).Value = "thailand"
$ws.Cells.Item(105, 3).Value = "thai-league-1"
$ws.Cells.Item(105, 4).Value = "2023-2024"
$ws.Cells.Item(105, 5).Value = 45276.54166666666
$ws.Cells.Item(105, 6).Value = "Uthai Thani"
$ws.Cells.Item(105, 7).Value = 2
$ws.Cells.Item(105, 8).Value = "Lamphun Warrior"
$ws.Cells.Item(105, 9).Value = 2
$ws.Cells.Item(105, 10).Value = 2.21
$ws.Cells.Item(105, 11).Value = "11/12/2023 11:12"
$ws.Cells.Item(105, 12).Value = 2.25
$ws.Cells.Item(105, 13).Value = "16/12/2023 12:52"
$ws.Cells.Item(105, 14).Value = 3.46
$ws.Cells.Item(105, 15).Value = "11/12/2023 11:12"
$ws.Cells.Item(105, 16).Value = 3.63
$ws.Cells.Item(105, 17).Value = "16/12/2023 12:52"
$ws.Cells.Item(105, 18).Value = 3.01
$ws.Cells.Item(105, 19).Value = "11/12/2023 11:12"
$ws.Cells.Item(105, 20).Value = 3.08
$ws.Cells.Item(105, 21).Value = "16/12/2023 12:52"
$ws.Cells.Item(105, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/uthai-thani-lamphun-warrior/KjvfGgfc/"

# Row 106
$ws.Cells.Item(106, 1).Value = 105
$ws.Cells.Item(106, 2).Value = "thailand"
$ws.Cells.Item(106, 3).Value = "thai-league-1"
$ws.Cells.Item(106, 4).Value = "2023-2024"
$ws.Cells.Item(106, 5).Value = 45276.58333333334
$ws.Cells.Item(106, 6).Value = "Nakhon Pathom"
$ws.Cells.Item(106, 7).Value = 1
$ws.Cells.Item(106, 8).Value = "Chiangrai Utd"
$ws.Cells.Item(106, 9).Value = 2
$ws.Cells.Item(106, 10).Value = 2.37
$ws.Cells.Item(106, 11).Value = "11/12/2023 14:12"
$ws.Cells.Item(106, 12).Value = 2.36
$ws.Cells.Item(106, 13).Value = "16/12/2023 13:56"
$ws.Cells.Item(106, 14).Value = 3.32
$ws.Cells.Item(106, 15).Value = "11/12/2023 14:12"
$ws.Cells.Item(106, 16).Value = 3.18
$ws.Cells.Item(106, 17).Value = "16/12/2023 13:56"
$ws.Cells.Item(106, 18).Value = 2.86
$ws.Cells.Item(106, 19).Value = "11/12/2023 14:12"
$ws.Cells.Item(106, 20).Value = 3.26
$ws.Cells.Item(106, 21).Value = "16/12/2023 13:56"
$ws.Cells.Item(106, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/nakhon-pathom-chiangrai-utd/6qOwE3DT/"

# Row 107
$ws.Cells.Item(107, 1).Value = 106
$ws.Cells.Item(107, 2).Value = "thailand"
$ws.Cells.Item(107, 3).Value = "thai-league-1"
$ws.Cells.Item(107, 4).Value = "2023-2024"
$ws.Cells.Item(107, 5).Value = 45277.47916666666
$ws.Cells.Item(107, 6).Value = "Police Tero"
$ws.Cells.Item(107, 7).Value = 2
$ws.Cells.Item(107, 8).Value = "Sukhothai"
$ws.Cells.Item(107, 9).Value = 3
$ws.Cells.Item(107, 10).Value = 2.1
$ws.Cells.Item(107, 11).Value = "12/12/2023 12:42"
$ws.Cells.Item(107, 12).Value = 2.49
$ws.Cells.Item(107, 13).Value = "17/12/2023 11:28"
$ws.Cells.Item(107, 14).Value = 3.53
$ws.Cells.Item(107, 15).Value = "12/12/2023 12:42"
$ws.Cells.Item(107, 16).Value = 3.74
$ws.Cells.Item(107, 17).Value = "17/12/2023 11:22"
$ws.Cells.Item(107, 18).Value = 3.17
$ws.Cells.Item(107, 19).Value = "12/12/2023 12:42"
$ws.Cells.Item(107, 20).Value = 2.66
$ws.Cells.Item(107, 21).Value = "17/12/2023 11:28"
$ws.Cells.Item(107, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/police-tero-sukhothai/AuX6DivG/"

# Row 108
$ws.Cells.Item(108, 1).Value = 107
$ws.Cells.Item(108, 2).Value = "thailand"
$ws.Cells.Item(108, 3).Value = "thai-league-1"
$ws.Cells.Item(108, 4).Value = "2023-2024"
$ws.Cells.Item(108, 5).Value = 45277.5
$ws.Cells.Item(108, 6).Value = "Port MTI FC"
$ws.Cells.Item(108, 7).Value = 4
$ws.Cells.Item(108, 8).Value = "Buriram"
$ws.Cells.Item(108, 9).Value = 1
$ws.Cells.Item(108, 10).Value = 2.57
$ws.Cells.Item(108, 11).Value = "11/12/2023 11:12"
$ws.Cells.Item(108, 12).Value = 2.54
$ws.Cells.Item(108, 13).Value = "17/12/2023 11:58"
$ws.Cells.Item(108, 14).Value = 3.37
$ws.Cells.Item(108, 15).Value = "11/12/2023 11:12"
$ws.Cells.Item(108, 16).Value = 3.8
$ws.Cells.Item(108, 17).Value = "17/12/2023 11:58"
$ws.Cells.Item(108, 18).Value = 2.57
$ws.Cells.Item(108, 19).Value = "11/12/2023 11:12"
$ws.Cells.Item(108, 20).Value = 2.57
$ws.Cells.Item(108, 21).Value = "17/12/2023 11:58"
$ws.Cells.Item(108, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/port-mti-fc-buriram-united-f-c/SGs2EXO9/"

# Row 109
$ws.Cells.Item(109, 1).Value = 108
$ws.Cells.Item(109, 2).Value = "thailand"
$ws.Cells.Item(109, 3).Value = "thai-league-1"
$ws.Cells.Item(109, 4).Value = "2023-2024"
$ws.Cells.Item(109, 5).Value = 45277.54166666666
$ws.Cells.Item(109, 6).Value = "Trat FC"
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = "Bangkok Utd"
$ws.Cells.Item(109, 9).Value = 2
$ws.Cells.Item(109, 10).Value = 5.73
$ws.Cells.Item(109, 11).Value = "11/12/2023 11:12"
$ws.Cells.Item(109, 12).Value = 5.58
$ws.Cells.Item(109, 13).Value = "17/12/2023 12:53"
$ws.Cells.Item(109, 14).Value = 4.48
$ws.Cells.Item(109, 15).Value = "11/12/2023 11:12"
$ws.Cells.Item(109, 16).Value = 4.26
$ws.Cells.Item(109, 17).Value = "17/12/2023 12:53"
$ws.Cells.Item(109, 18).Value = 1.45
$ws.Cells.Item(109, 19).Value = "11/12/2023 11:12"
$ws.Cells.Item(109, 20).Value = 1.57
$ws.Cells.Item(109, 21).Value = "17/12/2023 12:53"
$ws.Cells.Item(109, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/trat-fc-bangkok-utd/t0ujHZvi/"

